$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the old "filtre_visuels" column (column C). This shifts the
#    "instructions" column from D to C, matching the new 3-column layout.
$ws.Range("C1").EntireColumn.Delete()

# 2. Resize the (new) instructions column C to its final width and drop
#    the old best-fit sizing that used to live on column D.
$ws.Columns.Item(3).ColumnWidth = 94.75

# 3. Append the three new rows describing the "apps" json uniformisation,
#    one per filter (content rating, genres, type).
$ws.Cells.Item(12,1).Value = "jsons_train/unif_filtres/apps.json"
$ws.Cells.Item(12,2).Value = "jsons_train/unif_filtres/apps_same_as_content_rating.json"
$ws.Cells.Item(12,3).Value = "Uniformise le format de tous les filtres en te basant sur le format du filtre content rating"

$ws.Cells.Item(13,1).Value = "jsons_train/unif_filtres/apps.json"
$ws.Cells.Item(13,2).Value = "jsons_train/unif_filtres/apps_same_as_genre.json"
$ws.Cells.Item(13,3).Value = "Uniformise le format de tous les filtres en te basant sur le format du filtre genres"

$ws.Cells.Item(14,1).Value = "jsons_train/unif_filtres/apps.json"
$ws.Cells.Item(14,2).Value = "jsons_train/unif_filtres/apps_same_as_type.json"
$ws.Cells.Item(14,3).Value = "Uniformise le format de tous les filtres en te basant sur le format du filtre type"

# 4. The json path cell for the "browser_same_as_region" row ends up with
#    a distinct (but visually equivalent, General) number format in the
#    saved workbook; reproduce that extra style slot on B8.
$ws.Cells.Item(8,2).NumberFormatLocal = $ws.Cells.Item(8,2).NumberFormatLocal

# 5. Move the visible selection to C9, matching the saved cursor position.
$ws.Range("C9").Select() | Out-Null
